# Refresh the cryptocurrency price/volume columns (D = Price, E = Volume(1h))
# to match the latest scrape. Rows not listed below (e.g. 24 "Dai", 28
# "WrappedeETH") are unchanged. Price cells whose new text looks like a
# plain decimal number are forced back to Text via NumberFormat "@" so
# Excel does not silently reinterpret them as numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.627.65"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "2.570.52"
$ws.Range("E3").Value = "  -2.03%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.28"
$ws.Range("E5").Value = "  -1.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.21"
$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.98%  "

$ws.Range("D9").Value = "2.574.14"
$ws.Range("E9").Value = "  -1.90%  "

$ws.Range("E10").Value = "  +0.16%  "

$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  -1.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.16"
$ws.Range("E13").Value = "  -1.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.72"
$ws.Range("E14").Value = "  -3.13%  "

$ws.Range("D15").Value = "3.045.12"
$ws.Range("E15").Value = "  -2.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("E16").Value = "  -1.52%  "

$ws.Range("D17").Value = "66.555.27"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").Value = "2.579.32"
$ws.Range("E18").Value = "  -1.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.40"
$ws.Range("E19").Value = "  -6.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.72"
$ws.Range("E20").Value = "  -3.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.45"
$ws.Range("E21").Value = "  -1.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").Value = "  -1.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.59"
$ws.Range("E23").Value = "  -1.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.91"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.42"
$ws.Range("E26").Value = "  -0.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.92"
$ws.Range("E27").Value = "  -8.84%  "

$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "0.0₃0990"
$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.30"
$ws.Range("E31").Value = "  +5.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "528.98"
$ws.Range("E32").Value = "  -3.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.32"
$ws.Range("E33").Value = "  -2.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  -3.46%  "

$ws.Range("E35").Value = "  -3.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.46"
$ws.Range("E37").Value = "  -2.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.15"
$ws.Range("E38").Value = "  +0.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.75"
$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.358"
$ws.Range("E40").Value = "  -2.10%  "

$ws.Range("E41").Value = "  +2.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.77"
$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("E45").Value = "  +1.07%  "

$ws.Range("E46").Value = "  -3.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.92"
$ws.Range("E47").Value = "  -1.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.566"
$ws.Range("E48").Value = "  -2.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.71"
$ws.Range("E49").Value = "  -1.43%  "

$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0762"
$ws.Range("E51").Value = "  -1.18%  "
